$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F = FeltThermal, Column G = FeltMotion
$ws.Range("G2").Value = 0
$ws.Range("G10").Value = 1
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 1
$ws.Range("G51").Value = 1
$ws.Range("F63").Value = 1
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 1
$ws.Range("G91").Value = 1
$ws.Range("G119").Value = 1
$ws.Range("G132").Value = 1
$ws.Range("F145").Value = 1
$ws.Range("G148").Value = 1
$ws.Range("G163").Value = 1
$ws.Range("G177").Value = 1
